$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every "Wildtype" genotype value in column D (Genotype) with "Null"
for ($r = 2; $r -le 52; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Text -eq "Wildtype") {
        $cell.Value = "Null"
    }
}
